# Chiffres COVID-19 Valais - daily data update
# - Reset the frozen-pane scroll position back to the top of the data
#   (it had drifted down to row 252 while entering new rows).
# - Revise the SI-patient counts ("G") for 29 Sep - 19 Oct 2020 (rows 217-233):
#   a later correction added one more SI patient to most of those days,
#   while 19 Oct 2020 (row 233) was revised down by 4.
# - Correct/complete the daily entries for 4-12 Nov 2020 (rows 253-261):
#   new-case counts, intubated/SI-adjacent patient counts, SI patients and
#   hospital/extra-hospital death counts, including filling in row 261
#   which had not yet been entered.
#
# All of the totals/cumulative columns (B, H, J, K) are formulas and will
# recalculate automatically once the underlying inputs below are written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Scroll the frozen pane back to the top of the data (B3) ---
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 2

# --- SI patients ("Patients COVID-19 aux SI total"), column G ---
$ws.Cells.Item(217, 7).Value = 32
$ws.Cells.Item(218, 7).Value = 33
$ws.Cells.Item(219, 7).Value = 36
$ws.Cells.Item(220, 7).Value = 37
$ws.Cells.Item(221, 7).Value = 38
$ws.Cells.Item(222, 7).Value = 41
$ws.Cells.Item(223, 7).Value = 43
$ws.Cells.Item(224, 7).Value = 42
$ws.Cells.Item(225, 7).Value = 45
$ws.Cells.Item(226, 7).Value = 49
$ws.Cells.Item(227, 7).Value = 53
$ws.Cells.Item(228, 7).Value = 57
$ws.Cells.Item(229, 7).Value = 61
$ws.Cells.Item(230, 7).Value = 65
$ws.Cells.Item(231, 7).Value = 74
$ws.Cells.Item(232, 7).Value = 85
$ws.Cells.Item(233, 7).Value = 85

# --- Row 253 (04.11.2020): new positive cases ---
$ws.Range("C253").Value = 540

# --- Row 254 (05.11.2020): new positive cases ---
$ws.Range("C254").Value = 542

# --- Row 255 (06.11.2020): new positive cases ---
$ws.Range("C255").Value = 459

# --- Row 256 (07.11.2020): SI patients, new hospital deaths ---
$ws.Range("G256").Value = 294
$ws.Range("L256").Value = 5

# --- Row 257 (08.11.2020): new cases, intubated, hospitalised hors SI ---
$ws.Range("C257").Value = 191
$ws.Range("E257").Value = 37
$ws.Range("F257").Value = 23

# --- Row 258 (09.11.2020): new cases, intubated, hospitalised hors SI ---
$ws.Range("C258").Value = 554
$ws.Range("E258").Value = 34
$ws.Range("F258").Value = 25

# --- Row 259 (10.11.2020): new cases, intubated, hospitalised hors SI,
#     new hospital deaths, new extra-hospital deaths ---
$ws.Range("C259").Value = 426
$ws.Range("E259").Value = 36
$ws.Range("F259").Value = 27
$ws.Range("L259").Value = 3
$ws.Range("M259").Value = 6

# --- Row 260 (11.11.2020): new cases, intubated, hospitalised hors SI,
#     SI patients, new hospital deaths, new extra-hospital deaths ---
$ws.Range("C260").Value = 261
$ws.Range("E260").Value = 35
$ws.Range("F260").Value = 26
$ws.Range("G260").Value = 297
$ws.Range("L260").Value = 7
$ws.Range("M260").Value = 3

# --- Row 261 (12.11.2020): previously-empty row, now filled in ---
$ws.Range("C261").Value = 31
$ws.Range("E261").Value = 37
$ws.Range("F261").Value = 26
$ws.Range("G261").Value = 277
$ws.Range("L261").Value = 0
$ws.Range("M261").Value = 0
